$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "${GENERO}:" and "${NOMBRE_FUNCIONARIO}" runs switch from the
#    minorHAnsi theme font to an explicit "Century Gothic" ascii/hAnsi font
#    (bold / color / size / language stay untouched).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute('${GENERO}:', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    # Font.Name sets both w:ascii/w:hAnsi, but it also quietly resets the
    # complex-script font (w:cs) -- restore it to the original "Arial".
    $rng.Font.Name = "Century Gothic"
    $rng.Font.NameBi = "Arial"
}

$rng = $d.Content
$rng.Find.Execute('${NOMBRE_FUNCIONARIO}', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $rng.Font.Name = "Century Gothic"
    $rng.Font.NameBi = "Arial"
}

# ---------------------------------------------------------------------------
# 2) Rework "de Bs. ${MONTO}.- ${LITERAL}." into
#    "de Bs. ${MONTO}.- ${LITERAL} Bolivianos.", relocating the hidden
#    _GoBack bookmark so it wraps the new word "Bolivianos" instead of
#    "MONTO".
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")

# The two characters just before the bookmark are "${" (the start of the
# "${MONTO}" placeholder); rewrite that in place to read
# "${MONTO}.- ${LITERAL} " so the placeholder is whole and the following
# ".- ${LITERAL} " text is already present ahead of the bookmarked word.
$beforeRng = $d.Range($bm.Range.Start - 2, $bm.Range.Start)
$beforeRng.Text = '${MONTO}.- ${LITERAL} '

# Re-fetch the bookmark (its Range shifted because the text before it grew)
# and rename the word it wraps from "MONTO" to "Bolivianos".
$bm = $d.Bookmarks("_GoBack")
$bm.Range.Text = "Bolivianos"

# Finally, collapse the trailing "}.- ${LITERAL}." that used to follow
# "MONTO" down to a lone ".".
$bm = $d.Bookmarks("_GoBack")
$afterRng = $d.Range($bm.Range.End, $bm.Range.End + 15)
$afterRng.Text = "."
